# "Generate Report for Handoff" — a new handoff cycle has been kicked off for
# the two e2e markdown files. The first file was renamed/regenerated
# (c701babf... -> d1671c77...) and the second file too
# (eed7b301... -> ffff88617cf5...). Status flips from the old handback state
# to "Ready for handoff", the stale handback/target bookkeeping is cleared,
# and fresh handoff datetimes + xlf names are recorded.

$wb = $excel.ActiveWorkbook

$oldFile1 = "c701babf-ce22-4f8f-950f-b7fa872a7f71"
$newFile1 = "d1671c77-520e-4ab2-98a3-eed471feb362"
$oldFile2 = "eed7b301-1b14-4413-ba50-01a729aad8a4"
$newFile2 = "ffff88617cf5-5bee-48f6-a94a-84f11527f3e9"

$newStatus = "Ready for handoff"
$overviewDate = "2016-08-23 21:06:29"
$zhHandoffDate = "2016-08-23 21:06:24"
$deHandoffDate = "2016-08-23 21:06:29"
$epochDate = "0001-01-01 00:00:00"

$zhXlf1 = "$newFile1.82b700edb4aeb1bf4ff78bd011cb8d1ee4266f9a.zh-cn.xlf"
$deXlf1 = "$newFile1.82b700edb4aeb1bf4ff78bd011cb8d1ee4266f9a.de-de.xlf"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newFile1.md"
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Range("A3").Value = "$newFile2.md"
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = $overviewDate

# Re-point the two B-column hyperlinks (text + display) at the new file names.
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile1.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "e2e\$newFile2.md")

$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newFile1.md"
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("G2").Value = $zhXlf1
$wsZh.Range("H2").Value = $zhHandoffDate
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").ClearFormats()
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $epochDate

$wsZh.Range("A3").Value = "$newFile2.md"
$wsZh.Range("C3").Value = $newStatus
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhXlf1
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("I3").Value = ""
$wsZh.Range("I3").ClearFormats()
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $epochDate

# Drop the "Latest Target File" hyperlinks (I2/I3) entirely; keep A2/A3
# pointing at the source repo with the refreshed display text.
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md")

$wsZh.Columns.Item(3).ColumnWidth = 16.25
$wsZh.Columns.Item(9).ColumnWidth = 17.75
$wsZh.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newFile1.md"
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("G2").Value = $deXlf1
$wsDe.Range("H2").Value = $deHandoffDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").ClearFormats()
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $epochDate

$wsDe.Range("A3").Value = "$newFile2.md"
$wsDe.Range("C3").Value = $newStatus
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deXlf1
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("I3").ClearFormats()
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $epochDate

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile1.md", [Type]::Missing, [Type]::Missing, "$newFile1.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/903326a5ce7c39e30dac38cd4309fac06c5443bd/e2e/$newFile2.md", [Type]::Missing, [Type]::Missing, "$newFile2.md")

$wsDe.Columns.Item(3).ColumnWidth = 16.25
$wsDe.Columns.Item(9).ColumnWidth = 17.75
$wsDe.Columns.Item(10).ColumnWidth = 20.75
